$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Lucknow
$ws.Range("H2").Value = "88"
$ws.Range("K2").Value = "1.54"

# Row 3 - Bhopal
$ws.Range("B3").Value = "27.13"
$ws.Range("E3").Value = "Rain and moderate rain"
$ws.Range("H3").Value = "78"
$ws.Range("K3").Value = "4.87"

# Row 4 - Ajmer
$ws.Range("B4").Value = "28.99"
$ws.Range("E4").Value = "Clear and clear sky"
$ws.Range("H4").Value = "68"
$ws.Range("K4").Value = "8.03"

# Row 5 - Coimbatore
$ws.Range("B5").Value = "25.00"
$ws.Range("H5").Value = "94"
$ws.Range("K5").Value = "4.60"

# Row 7 - Kolkata
$ws.Range("K7").Value = "2.10"
